$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H58").Value = 1383.2858
$ws.Range("I58").Value = 954.75
$ws.Range("J58").Value = 1554.7
$ws.Range("K58").Value = 2864.25
$ws.Range("L58").Value = 4664.1
$ws.Range("M58").Value = -2714.25
$ws.Range("N58").Value = -4964.1

$ws.Range("H62").Value = 7550.0625
$ws.Range("I62").Value = 5399.2856
$ws.Range("K62").Value = 5399.2856
$ws.Range("M62").Value = -4775.2856

$ws.Range("H65").Value = 7550.0625
$ws.Range("I65").Value = 5399.2856
$ws.Range("K65").Value = 26996.428
$ws.Range("M65").Value = -23876.428

$ws.Range("H69").Value = 8699.375
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8699.375
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 26098.125
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -27846.125

$ws.Range("H72").Value = 8699.375
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8699.375
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 78294.375
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -87030.375

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 2910.3845
$ws.Range("I107").Value = 2092.2222
$ws.Range("J107").Value = 4751.25
$ws.Range("K107").Value = 2092.2222
$ws.Range("L107").Value = 4751.25
$ws.Range("M107").Value = -172.2222000000002
$ws.Range("N107").Value = -8591.25

$ws.Range("H125").Value = 2986.4
$ws.Range("I125").Value = 2644
$ws.Range("J125").Value = 3500
$ws.Range("K125").Value = 23796
$ws.Range("L125").Value = 31500
$ws.Range("M125").Value = -21336
$ws.Range("N125").Value = -36420

$ws.Range("H132").Value = 3115.4666
$ws.Range("I132").Value = 1927.0714
$ws.Range("K132").Value = 5781.2142
$ws.Range("M132").Value = -3251.2142

$ws.Range("H135").Value = 1774.6364
$ws.Range("I135").Value = 1792.5714
$ws.Range("J135").Value = 1743.25
$ws.Range("K135").Value = 16133.1426
$ws.Range("L135").Value = 15689.25
$ws.Range("M135").Value = -13598.1426
$ws.Range("N135").Value = -20759.25

$ws.Range("H138").Value = 3524.3333
$ws.Range("I138").Value = 3161.5833
$ws.Range("J138").Value = 3887.0833
$ws.Range("K138").Value = 9484.749899999999
$ws.Range("L138").Value = 11661.2499
$ws.Range("M138").Value = -4344.749899999999
$ws.Range("N138").Value = -21941.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5805.387
$ws.Range("J61").Value = 15750
$ws.Range("L61").Value = 15750
$ws.Range("N61").Value = -16174

$ws.Range("H86").Value = 67500
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 67500
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H132").Value = 7071.2666
$ws.Range("I132").Value = 5219.72
$ws.Range("K132").Value = 15659.16
$ws.Range("M132").Value = -13129.16

$ws.Range("H136").Value = 5805.387
$ws.Range("J136").Value = 15750
$ws.Range("L136").Value = 47250
$ws.Range("N136").Value = -52350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 938.1539
$ws.Range("I94").Value = 955.2222
$ws.Range("J94").Value = 899.75
$ws.Range("K94").Value = 955.2222
$ws.Range("L94").Value = 899.75
$ws.Range("M94").Value = -504.2222
$ws.Range("N94").Value = -1801.75

$ws.Range("H105").Value = 19963.215
$ws.Range("I105").Value = 22504.8
$ws.Range("K105").Value = 22504.8
$ws.Range("M105").Value = -20757.8

$ws.Range("H134").Value = 2714.8928
$ws.Range("I134").Value = 1625.125
$ws.Range("K134").Value = 4875.375
$ws.Range("M134").Value = -2340.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1375.6875
$ws.Range("I22").Value = 820.1539
$ws.Range("J22").Value = 3783
$ws.Range("K22").Value = 820.1539
$ws.Range("L22").Value = 3783
$ws.Range("M22").Value = -470.1539
$ws.Range("N22").Value = -4483

$ws.Range("H31").Value = 27029.957
$ws.Range("I31").Value = 4309
$ws.Range("K31").Value = 4309
$ws.Range("M31").Value = -4014

$ws.Range("H34").Value = 27029.957
$ws.Range("I34").Value = 4309
$ws.Range("K34").Value = 4309
$ws.Range("M34").Value = -4107

$ws.Range("H58").Value = 4281.6787
$ws.Range("I58").Value = 1963.2
$ws.Range("K58").Value = 1963.2
$ws.Range("M58").Value = -1760.2

$ws.Range("H95").Value = 21234.715
$ws.Range("J95").Value = 21234.715
$ws.Range("L95").Value = 21234.715
$ws.Range("N95").Value = -26726.715

$ws.Range("H136").Value = 4281.6787
$ws.Range("I136").Value = 1963.2
$ws.Range("K136").Value = 5889.6
$ws.Range("M136").Value = -3339.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1489.4
$ws.Range("J52").Value = 1489.4
$ws.Range("L52").Value = 4468.200000000001
$ws.Range("N52").Value = -5000.200000000001

$ws.Range("H68").Value = 2795.9583
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2874.0435
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 8622.130500000001
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -10244.1305

$ws.Range("H71").Value = 2795.9583
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2874.0435
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 25866.3915
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -33978.3915

$ws.Range("H108").Value = 4460.909
$ws.Range("I108").Value = 1749.375
$ws.Range("K108").Value = 5248.125
$ws.Range("M108").Value = -2368.125

$ws.Range("H111").Value = 2841.3333
$ws.Range("I111").Value = 2841.3333
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8523.999899999999
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5456.999899999999
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1989
$ws.Range("I107").Value = 722.375
$ws.Range("K107").Value = 722.375
$ws.Range("M107").Value = 1197.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1854.625
$ws.Range("I16").Value = 1854.625
$ws.Range("K16").Value = 1854.625
$ws.Range("M16").Value = -1684.625

$ws.Range("H40").Value = 8145.1724
$ws.Range("I40").Value = 7843.143
$ws.Range("J40").Value = 8938
$ws.Range("K40").Value = 7843.143
$ws.Range("L40").Value = 8938
$ws.Range("M40").Value = -7707.143
$ws.Range("N40").Value = -9210

$ws.Range("H61").Value = 4035.818
$ws.Range("I61").Value = 4035.818
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4035.818
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3833.818
$ws.Range("N61").ClearContents()

$ws.Range("H88").Value = 186666.67
$ws.Range("I88").Value = 60000
$ws.Range("J88").Value = 250000
$ws.Range("K88").Value = 60000
$ws.Range("L88").Value = 250000
$ws.Range("M88").Value = -59572
$ws.Range("N88").Value = -250856

$ws.Range("H91").Value = 186666.67
$ws.Range("I91").Value = 60000
$ws.Range("J91").Value = 250000
$ws.Range("K91").Value = 60000
$ws.Range("L91").Value = 250000
$ws.Range("M91").Value = -58518
$ws.Range("N91").Value = -252964

$ws.Range("H93").Value = 15559.483
$ws.Range("I93").Value = 12889.167
$ws.Range("K93").Value = 12889.167
$ws.Range("M93").Value = -11641.167

$ws.Range("H100").Value = 2450
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -359
$ws.Range("N100").Value = -5082

$ws.Range("H113").Value = 4035.818
$ws.Range("I113").Value = 4035.818
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4035.818
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1865.818
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 37791.89
$ws.Range("I122").Value = 43179.145
$ws.Range("J122").Value = 18936.5
$ws.Range("K122").Value = 129537.435
$ws.Range("L122").Value = 56809.5
$ws.Range("M122").Value = -127087.435
$ws.Range("N122").Value = -61709.5

$ws.Range("H136").Value = 6723.6284
$ws.Range("I136").Value = 4145.3887
$ws.Range("J136").Value = 9453.529
$ws.Range("K136").Value = 12436.1661
$ws.Range("L136").Value = 28360.587
$ws.Range("M136").Value = -9886.166100000002
$ws.Range("N136").Value = -33460.587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 6348
$ws.Range("I55").Value = 996.3333
$ws.Range("J55").Value = 11699.667
$ws.Range("K55").Value = 996.3333
$ws.Range("L55").Value = 11699.667
$ws.Range("M55").Value = -719.3333
$ws.Range("N55").Value = -12253.667

$ws.Range("H81").Value = 11950.125
$ws.Range("I81").Value = 7859.4
$ws.Range("K81").Value = 15718.8
$ws.Range("M81").Value = -14657.8

$ws.Range("H84").Value = 11950.125
$ws.Range("I84").Value = 7859.4
$ws.Range("K84").Value = 78594
$ws.Range("M84").Value = -73290

$ws.Range("H132").Value = 6216.5884
$ws.Range("I132").Value = 3150.3103
$ws.Range("K132").Value = 9450.930899999999
$ws.Range("M132").Value = -6920.930899999999
